$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G4").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G14").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G18").Value = "Shivamogga (Shimoga)"
$ws.Range("G29").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G30").Value = "Shivamogga (Shimoga)"
$ws.Range("G38").Value = "Shivamogga (Shimoga)"
$ws.Range("G39").Value = "Shivamogga (Shimoga)"
$ws.Range("G42").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G43").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G45").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G53").Value = "Chikkamagaluru (Chikmagalur)"
